$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(150).Insert()

$ws.Range("A150").Value = 8
$ws.Range("B150").Value = "Terminal La Palmera de La Serena"
$ws.Range("C150").Value = "Coquimbo"
$ws.Range("D150").Value = 44985
$ws.Range("E150").Value = 4
$ws.Range("F150").Value = 100112044
$ws.Range("G150").Value = "Perejil"
$ws.Range("H150").Value = "Sin especificar"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 2500
$ws.Range("K150").Value = 2000
$ws.Range("L150").Value = 2500
$ws.Range("M150").Value = 2250
$ws.Range("N150").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O150").Value = "Provincia del Elquí"
$ws.Range("P150").Value = 1500
$ws.Range("Q150").Value = 1.5
$ws.Range("R150").Value = "Hortaliza"

Write-Output "done"
